$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '44.136.91'
$ws.Cells.Item(2, 5).Value = '  +2.50%  '

$ws.Cells.Item(3, 4).Value = '2.282.41'
$ws.Cells.Item(3, 5).Value = '  +2.57%  '

$ws.Cells.Item(4, 5).Value = '  -0.14%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '318.95'
$ws.Cells.Item(5, 5).Value = '  +0.15%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '103.15'
$ws.Cells.Item(6, 5).Value = '  +4.29%  '

$ws.Cells.Item(7, 5).Value = '  +1.12%  '

$ws.Cells.Item(8, 5).Value = '  -0.16%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.572'
$ws.Cells.Item(9, 5).Value = '  +1.30%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '38.67'
$ws.Cells.Item(10, 5).Value = '  +5.89%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0839'
$ws.Cells.Item(11, 5).Value = '  +1.57%  '

$ws.Cells.Item(12, 5).Value = '  +2.10%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.107'
$ws.Cells.Item(13, 5).Value = '  +2.07%  '

$ws.Cells.Item(14, 4).Value = '2.630.56'
$ws.Cells.Item(14, 5).Value = '  +2.55%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.876'
$ws.Cells.Item(15, 5).Value = '  +1.10%  '

$ws.Cells.Item(16, 5).Value = '  +3.80%  '

$ws.Cells.Item(17, 4).Value = '2.289.02'
$ws.Cells.Item(17, 5).Value = '  +2.81%  '

$ws.Cells.Item(18, 4).Value = '44.058.52'
$ws.Cells.Item(18, 5).Value = '  +2.79%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '14.45'
$ws.Cells.Item(19, 5).Value = '  -0.98%  '

$ws.Cells.Item(20, 5).Value = '  +3.48%  '

$ws.Cells.Item(21, 5).Value = '  +3.26%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '66.28'
$ws.Cells.Item(22, 5).Value = '  +1.69%  '

$ws.Cells.Item(23, 5).Value = '  +1.40%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '238.34'
$ws.Cells.Item(24, 5).Value = '  +1.00%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.20'
$ws.Cells.Item(25, 5).Value = '  +3.69%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.29%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.31'
$ws.Cells.Item(27, 5).Value = '  +1.40%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '39.38'
$ws.Cells.Item(28, 5).Value = '  +16.39%  '

$ws.Cells.Item(29, 5).Value = '  +1.37%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.52'
$ws.Cells.Item(30, 5).Value = '  +3.64%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '162.31'
$ws.Cells.Item(31, 5).Value = '  +3.21%  '

$ws.Cells.Item(32, 5).Value = '  +0.29%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0880'
$ws.Cells.Item(33, 5).Value = '  -0.35%  '

$ws.Cells.Item(34, 5).Value = '  -1.69%  '

$ws.Cells.Item(35, 5).Value = '  -1.02%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.06'
$ws.Cells.Item(36, 5).Value = '  +1.64%  '

$ws.Cells.Item(37, 5).Value = '  -1.06%  '

$ws.Cells.Item(38, 5).Value = '  +1.39%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.109'
$ws.Cells.Item(39, 5).Value = '  +4.66%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.88'
$ws.Cells.Item(40, 5).Value = '  +5.97%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '15.73'
$ws.Cells.Item(41, 5).Value = '  +29.59%  '

$ws.Cells.Item(42, 5).Value = '  +0.77%  '

$ws.Cells.Item(43, 5).Value = '  -0.01%  '

$ws.Cells.Item(44, 4).Value = '1.776.41'
$ws.Cells.Item(44, 5).Value = '  -4.43%  '

$ws.Cells.Item(45, 5).Value = '  +0.73%  '

$ws.Cells.Item(46, 2).Value = 'THORChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '5.43'
$ws.Cells.Item(46, 5).Value = '  -1.22%  '

$ws.Cells.Item(47, 2).Value = 'BitcoinSV'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '85.16'
$ws.Cells.Item(47, 5).Value = '  -3.35%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '8.95'
$ws.Cells.Item(48, 5).Value = '  +3.77%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '59.62'
$ws.Cells.Item(49, 5).Value = '  -1.54%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '74.81'
$ws.Cells.Item(50, 5).Value = '  -4.78%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '104.75'
$ws.Cells.Item(51, 5).Value = '  +4.13%  '
